# Update the lattice-multiplication exercise table: every cell in the
# 5x3 table gets new multiplication problems (and matching partial-product
# lattice digits), while the surrounding formatting/structure is unchanged.
#
# Each cell's text is five lines (separated by Word's manual line-break,
# vertical-tab / char 11, i.e. backtick-v) :
#   "<a> x <b>"
#   "  <b digit> <b digit>"   (right-aligned under the two-digit factor)
#   "  ----"
#   "<d>|    |"
#   "<d>|    |"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-LatticeCell {
    param($row, $col, $text)
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $text
}

Set-LatticeCell 1 1 "47 x 94`v  9    4`v  ----`v4|    |`v7|    |"
Set-LatticeCell 1 2 "88 x 84`v  8    4`v  ----`v8|    |`v8|    |"
Set-LatticeCell 1 3 "76 x 50`v  5    0`v  ----`v7|    |`v6|    |"

Set-LatticeCell 2 1 "38 x 92`v  9    2`v  ----`v3|    |`v8|    |"
Set-LatticeCell 2 2 "93 x 37`v  3    7`v  ----`v9|    |`v3|    |"
Set-LatticeCell 2 3 "44 x 60`v  6    0`v  ----`v4|    |`v4|    |"

Set-LatticeCell 3 1 "11 x 29`v  2    9`v  ----`v1|    |`v1|    |"
Set-LatticeCell 3 2 "55 x 30`v  3    0`v  ----`v5|    |`v5|    |"
Set-LatticeCell 3 3 "66 x 76`v  7    6`v  ----`v6|    |`v6|    |"

Set-LatticeCell 4 1 "80 x 72`v  7    2`v  ----`v8|    |`v0|    |"
Set-LatticeCell 4 2 "87 x 56`v  5    6`v  ----`v8|    |`v7|    |"
Set-LatticeCell 4 3 "28 x 21`v  2    1`v  ----`v2|    |`v8|    |"

Set-LatticeCell 5 1 "51 x 66`v  6    6`v  ----`v5|    |`v1|    |"
Set-LatticeCell 5 2 "12 x 33`v  3    3`v  ----`v1|    |`v2|    |"
Set-LatticeCell 5 3 "69 x 71`v  7    1`v  ----`v6|    |`v9|    |"

Write-Output "Lattice multiplication exercises updated"
